# Tugas UTS-Agile.docx edit
# ------------------------
# In the paragraph "Deskripsi singkat tengtang aplikasi e-Doc":
#   - remove the eight stray <w:proofErr w:type="spellStart|spellEnd"/>
#     markers that bracket its runs (they are Word's "needs spell-check"
#     bookmarks, carry no visible content, and are not reachable through
#     the normal text/Find object model)
#   - fix the trailing text " e-Doc" -> " e-Docter"
#
# Because proofErr markers aren't exposed as editable content, the paragraph
# is rebuilt in one shot with Range.InsertXML using an exact copy of its
# original OOXML (same runs, same run/paragraph formatting) minus the
# proofErr markers and with the corrected trailing text.

$d = $word.ActiveDocument

$targetParagraphXml = "<?xml version=`"1.0`" standalone=`"yes`"?><?mso-application progid=`"Word.Document`"?><pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`"><w:body><w:p w:rsidR=`"4F653FCB`" w:rsidP=`"70A4EC72`" w:rsidRDefault=`"4F653FCB`" w14:paraId=`"4EFD5879`" w14:textId=`"75A92BEC`"><w:pPr><w:spacing w:line=`"257`" w:lineRule=`"auto`" /><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr></w:pPr><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t>Deskripsi</w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t>singkat</w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t>tengtang</w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t xml:space=`"preserve`"> </w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t>aplikasi</w:t></w:r><w:r w:rsidRPr=`"70A4EC72`"><w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`" w:eastAsia=`"Times New Roman`" w:cs=`"Times New Roman`" /><w:sz w:val=`"28`" /><w:szCs w:val=`"28`" /></w:rPr><w:t xml:space=`"preserve`"> e-Docter</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text
    if ($paraText -like "*Deskripsi*singkat*tengtang*aplikasi*e-Doc*") {
        $para.Range.InsertXML($targetParagraphXml)
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the 'Deskripsi singkat tengtang aplikasi e-Doc' paragraph"
}
